$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Paragraph 1 ("School of Computer ... Cardiff University"):
# the heading text used to be split across two runs with a stray, empty
# "_GoBack" bookmark sitting between them. Drop the leftover bookmark and
# normalize the heading back into a single run of text.
# ---------------------------------------------------------------------------
$bmGoBack = $d.Bookmarks.Item("_GoBack")
$bmGoBack.Delete()

$heading = $d.Paragraphs.Item(1).Range
$heading.Find.Execute("School of Computer Science & Informatics, Cardiff University", $true, $false, $false, $false, $false, $true, 1, $false, "School of Computer Science & Informatics, Cardiff University", 2)

# ---------------------------------------------------------------------------
# Paragraph 2 ("Employer's SFIA Assessment - Professional IT Skills - Animation
# Development"): remove the erroneous "Animation Development" suffix (and the
# dash that introduced it) along with the now-unused "_Hlk29918752" bookmark,
# while keeping the still-referenced "_Hlk29918582" bookmark in place at the
# end of the (now shorter) title.
# ---------------------------------------------------------------------------
$bmStray = $d.Bookmarks.Item("_Hlk29918752")
$bmStray.Delete()

$bmKeep = $d.Bookmarks.Item("_Hlk29918582")
$titlePara = $d.Paragraphs.Item(2).Range

# Remove "Animation Development" (the text right after the bookmark, up to
# - but not including - the paragraph mark at the end of the title).
$afterBookmark = $d.Range($bmKeep.Start, $titlePara.End - 1)
$afterBookmark.Text = ""

# Remove the trailing " - " that used to lead into "Animation Development".
$beforeBookmark = $d.Range($bmKeep.Start - 3, $bmKeep.Start)
$beforeBookmark.Text = ""
